$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"
$wsOverview.Range("D2").Value = "2016-47-19 20:47:47"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"
$wsZhCn.Range("D2").Value = "1d174f44-0cbd-470f-adef-d1d62a5a2945.45fd3358a5866ca9c68b71f3245768bcb6581ccd.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-19 20:47:44"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"
$wsDeDe.Range("D2").Value = "1d174f44-0cbd-470f-adef-d1d62a5a2945.45fd3358a5866ca9c68b71f3245768bcb6581ccd.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-19 20:47:47"

# --- Update hyperlink display text to match new file names ---
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"

$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"
$wsZhCn.Hyperlinks.Item(3).TextToDisplay = "1d174f44-0cbd-470f-adef-d1d62a5a2945.45fd3358a5866ca9c68b71f3245768bcb6581ccd.zh-cn.xlf"

$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "1d174f44-0cbd-470f-adef-d1d62a5a2945.md"
$wsDeDe.Hyperlinks.Item(3).TextToDisplay = "1d174f44-0cbd-470f-adef-d1d62a5a2945.45fd3358a5866ca9c68b71f3245768bcb6581ccd.de-de.xlf"

$wb.Save()
